# Add sequential requirement IDs (e.g. "psm-1.1", "psm-1.2", ...) to column A
# of the "FUNC Reqs" sheet. Each block of requirement rows sits directly below
# a bold section-header row (column B holds the section title, e.g.
# "1.  Capability to conduct identity verification"); within a section the
# rows are numbered 1..N in "psm-<section>.<n>" form.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("FUNC Reqs")

# Row numbers of the section-header rows (col B holds the section title,
# col A is blank on these rows).
$sectionHeaderRows = @(2, 8, 28, 53, 70, 75, 78)

# Last data row on the sheet (row 92 is the final requirement row).
$lastDataRow = 92

for ($i = 0; $i -lt $sectionHeaderRows.Count; $i++) {
    $headerRow = $sectionHeaderRows[$i]
    $sectionNumber = $i + 1

    $startRow = $headerRow + 1
    if ($i + 1 -lt $sectionHeaderRows.Count) {
        $endRow = $sectionHeaderRows[$i + 1] - 1
    } else {
        $endRow = $lastDataRow
    }

    $reqNumber = 1
    for ($row = $startRow; $row -le $endRow; $row++) {
        $ws.Cells.Item($row, 1).Value = "psm-$sectionNumber.$reqNumber"
        $reqNumber = $reqNumber + 1
    }
}
